$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column L (nemad / ticker symbol) currently repeats the full company name
# "نفت سپاهان" for every data row (2-45). Update it to the ticker symbol
# "شسپا" while leaving column M (company_name) untouched.
for ($r = 2; $r -le 45; $r++) {
    $ws.Range("L$r").Value = "شسپا"
}

# Also select L7, matching the workbook's recorded active cell/selection.
$ws.Range("L7").Select()
